# Sprint 3 Backlog: add three new backlog rows (15-17) to the Product Backlog
# sheet, and move the active selection to E2, matching the author's commit
# "Updated all documentation, and added Sprint 3 Backlog".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: Game / Player story about finishing the game with online highscores.
$ws.Cells.Item(15, 1).Value = "Game"
$ws.Cells.Item(15, 2).Value = "Player"
$ws.Cells.Item(15, 3).Value = "have a complete game"
$ws.Cells.Item(15, 4).Value = "I can finally play"
$ws.Cells.Item(15, 5).Value = "with online highscores"
$ws.Cells.Item(15, 6).Value = "High"
$ws.Cells.Item(15, 7).Value = "Open"

# Row 16: Non-functional / Developer story about clean, reviewable code.
$ws.Cells.Item(16, 1).Value = "Non-functional"
$ws.Cells.Item(16, 2).Value = "Developer"
$ws.Cells.Item(16, 3).Value = "have clean and tidy final code"
$ws.Cells.Item(16, 4).Value = "my work may be reviewed favorably"
$ws.Cells.Item(16, 6).Value = "High"
$ws.Cells.Item(16, 7).Value = "Open"

# Row 17: Game / Player story about a shared online scoreboard.
$ws.Cells.Item(17, 1).Value = "Game"
$ws.Cells.Item(17, 2).Value = "Player"
$ws.Cells.Item(17, 3).Value = "have all my fellow players share a scoreboard with me on a server"
$ws.Cells.Item(17, 4).Value = "I may compare my skills "
$ws.Cells.Item(17, 6).Value = "High"
$ws.Cells.Item(17, 7).Value = "Open"

# Move the active selection, as it was left after the edit session.
$ws.Range("E2").Select() | Out-Null
